# issue #5: add legislator_id, name, date into dataframe
#
# The "股票" (stocks) worksheet gets three new trailing columns:
#   H = date             (text, e.g. "2013-11-01")
#   I = legislator_name   (text, e.g. "鄭汝芬")
#   J = legislator_id     (number, e.g. 1713)
# for every existing data row.

$wb = $excel.ActiveWorkbook

# Locate the "股票" worksheet by name (falls back to the 5th sheet,
# which is where it lives in this workbook, if the name can't be found).
$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "股票") {
        $ws = $sheet
    }
}
if ($ws -eq $null) {
    $ws = $wb.Worksheets.Item(5)
}

$lastRow = 13
$dateCol = 8
$nameCol = 9
$idCol = 10

$dateValue = "2013-11-01"
$legislatorName = "鄭汝芬"
$legislatorId = 1713

# --- Header row (row 1): copy the bold / centered / bordered look of the
# existing header cells (style index 1) onto the new header cells. ---
$h1 = $ws.Cells.Item(1, $dateCol)
$h1.Borders.LineStyle = 1
$h1.Borders.Weight = 2
$h1.HorizontalAlignment = -4108
$h1.VerticalAlignment = -4160
$h1.Font.Bold = $true
$h1.Copy()
$ws.Range($ws.Cells.Item(1, $nameCol), $ws.Cells.Item(1, $idCol)).PasteSpecial(-4122)

$ws.Cells.Item(1, $dateCol).Value = "date"
$ws.Cells.Item(1, $nameCol).Value = "legislator_name"
$ws.Cells.Item(1, $idCol).Value = "legislator_id"

# Store the date column as text so "2013-11-01" is kept as a literal
# string instead of being re-interpreted as a date serial number.
$ws.Range($ws.Cells.Item(2, $dateCol), $ws.Cells.Item($lastRow, $dateCol)).NumberFormat = "@"

# --- Data rows: fill in the new values for every existing record. ---
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $dateCol).Value = $dateValue
    $ws.Cells.Item($r, $nameCol).Value = $legislatorName
    $ws.Cells.Item($r, $idCol).Value = $legislatorId
}
